$d = $word.ActiveDocument

# --- 1. Remove the "Date" paragraph ("27 marca 2018") entirely ---
$removed = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.ParagraphStyle.NameLocal -eq "Date") {
        $p.Range.Delete()
        $removed = $true
        break
    }
}
Write-Output "Date paragraph removed: $removed"

# --- 2. Replace "Aim of the study" paragraph text ---
$r2 = $d.Content.Find.Execute(
    "The aim of the study was to investigate prevalence and factors associated with illegal substance use among Polish university students and to compare the extent of the phenomenon between different faculties.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The aim of the study was to investigate the prevalence and factors associated with illegal substance use among Polish university students.",
    2
)
Write-Output "Aim paragraph replaced: $r2"

# --- 3. Replace "Material and methods" paragraph text ---
$r3 = $d.Content.Find.Execute(
    "An internet questionnaire was distributed through social media. Data on age, gender, year of study, faculty, and use of certain psychoactive substances (both legal and illegal) during studying were collected. Respondents were able to add custom responses to the multiple-choice question about the substances used (thus mentioned drugs were then manually classified as legal or illegal). Additionally, two 5-point Likert-like scale based questions were introduced in order to quantify actual interest in the subject of their studies and potential social pressure to choose certain faculty or to take up higher education in general. The Shapiro-Wilk W test was used to test the normality of distribution; Mann-Whitney U and Chi-Square were used where appropriate. Logistic regression was used to determine the factors associated with illegal substance use. Analyses were performed using R v. 3.4.3. Database and scripts are available upon request.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data on age, gender, year of study, faculty, and use of certain psychoactive substances during studying were collected via an online questionnaire. Respondents were able to list custom substances (thus mentioned drugs were then manually classified as legal or illegal). Additionally, two 5-point Likert-like scale based questions were introduced in order to quantify actual interest in the subject of their studies and potential social pressure to choose certain faculty or to take up higher education in general. Mann-Whitney U and Chi-Square tests were used where appropriate. Logistic regression was used to determine the factors associated with illegal substance use. Database and scripts (R v. 3.4.3) are available upon request.",
    2
)
Write-Output "Material and methods paragraph replaced: $r3"

# --- 4. Replace "Results" paragraph text ---
$r4 = $d.Content.Find.Execute(
    "In total, 792 university students from over a hundred different faculties responded to the survey, medical students comprised 34.34% (N = 272) of the sample. In general, students who claimed to have used illegal substances during studying had experienced more pressure to take up their studies than those who claimed otherwise (M = 2.34 ± 1.52 vs. M = 1.83 ± 1.17, p=0.016). Additionally, they were less interested in the subject they were studying (3.68 ± 1.2 vs. 4.06 ± 0.9, p= 0.039). There were no significant differences in the age or year of study between these two groups (p= 0.821 and p=0.203, respectively). A smaller percentage of medical students, in comparison to other students, claimed to have used illegal substances, however the difference was not statistically significant (5.51% vs. 7.78%, p=0.276). Four independent factors associated with illegal substance use were identified by means of logistic regression: male sex (OR 2.13, 95% CI 1.18-3.77), caffeine consumption during studying (OR 3.72, 95% CI 1.67- 9.94), higher degree of being pressured to study at a certain faculty (OR 1.30, 95% CI 1.05-1.60) and not studying at a medical university (OR 2.38, 95% CI 1.32-4.35).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In total, 792 university students responded to the survey, medical students comprised 34.34% (N = 272) of the sample. Students who claimed to have used illegal substances during studying had experienced more pressure to take up their studies than those who claimed otherwise (M = 2.34 ± 1.52 vs. M = 1.83 ± 1.17, p=0.016). Additionally, they were less interested in the subject they were studying (3.68 ± 1.2 vs. 4.06 ± 0.9, p= 0.039). A smaller percentage of medical students, in comparison to other students, claimed to have used illegal substances, however the difference was not statistically significant (5.51% vs. 7.78%, p=0.276). Four independent factors associated with illegal substance use were identified by means of logistic regression: male sex (OR 2.13, 95% CI 1.18-3.77), caffeine consumption during studying (OR 3.72, 95% CI 1.67- 9.94), higher degree of pressure to study at a certain faculty (OR 1.30, 95% CI 1.05-1.60) and not studying at a medical university (OR 2.38, 95% CI 1.32-4.35).",
    2
)
Write-Output "Results paragraph replaced: $r4"
